# fixing Pembro and Atezo params
# - add missing "1/d" units label for the k12 / k21 transit-rate rows
# - update keL (row 27) and keTL (row 28) parameter values; the F26/F30
#   formulas that depend on them (ksynL and Kss_TL) recalc automatically
# - move the active selection to F29

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Units column ("1/d") was missing for rows 9 (k12) and 10 (k21)
$ws.Range("G9").Value = "1/d"
$ws.Range("G10").Value = "1/d"

# keL: 20 -> 6
$ws.Range("F27").Value = 6

# keTL: 20 -> 8
$ws.Range("F28").Value = 8

# Update the sheet's active cell/selection to match the saved view state
[void]$ws.Range("F29").Select()
